$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 4.33
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 1.8
$ws.Range("J4").Value = 4.75
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 2.4
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("W4").Value = 13
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 15
$ws.Range("AA4").Value = 34
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 7
$ws.Range("AI4").Value = 8.5
$ws.Range("AJ4").Value = 8.5
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 15
$ws.Range("AM4").Value = 26
$ws.Range("AN4").Value = 6
$ws.Range("AO4").Value = 23
$ws.Range("AQ4").Value = 81
$ws.Range("AR4").Value = 101
$ws.Range("AX4").Value = 3.75
$ws.Range("AY4").Value = 9.5
$ws.Range("BA4").Value = 34
$ws.Range("G5").Value = 2.15
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.6
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 1.95
$ws.Range("V5").Value = 1.8
$ws.Range("X5").Value = 9.5
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 19
$ws.Range("AA5").Value = 19
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 6
$ws.Range("AG5").Value = 351
$ws.Range("AH5").Value = 9.5
$ws.Range("AI5").Value = 17
$ws.Range("AL5").Value = 29
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 23
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 67
$ws.Range("BC5").Value = 251
$ws.Range("G6").Value = 1.36
$ws.Range("I6").Value = 9.5
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 3
$ws.Range("S6").Value = 1.5
$ws.Range("T6").Value = 2.5
$ws.Range("U6").Value = 2.75
$ws.Range("V6").Value = 1.4
$ws.Range("W6").Value = 4.75
$ws.Range("X6").Value = 5
$ws.Range("Z6").Value = 8
$ws.Range("AC6").Value = 7.5
$ws.Range("AD6").Value = 9.5
$ws.Range("AE6").Value = 34
$ws.Range("AF6").Value = 151
$ws.Range("AJ6").Value = 29
$ws.Range("AK6").Value = 151
$ws.Range("AM6").Value = 101
$ws.Range("AN6").Value = 3
$ws.Range("AO6").Value = 6.5
$ws.Range("AP6").Value = 26
$ws.Range("AS6").Value = 301
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 12
$ws.Range("AX6").Value = 10
$ws.Range("BA6").Value = 351
$ws.Range("BB6").Value = 401
$ws.Range("G9").Value = 1.57
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.2
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 6.5
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 11
$ws.Range("AA9").Value = 15
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 7.5
$ws.Range("AH9").Value = 13
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 67
$ws.Range("AL9").Value = 51
$ws.Range("AN9").Value = 3.4
$ws.Range("AO9").Value = 8
$ws.Range("AP9").Value = 21
$ws.Range("AQ9").Value = 26
$ws.Range("AR9").Value = 51
$ws.Range("AT9").Value = 2.63
$ws.Range("AX9").Value = 7.5
$ws.Range("AY9").Value = 34
$ws.Range("BA9").Value = 126
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 2.88
$ws.Range("I10").Value = 4.5
$ws.Range("J10").Value = 3
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 1.17
$ws.Range("N10").Value = 5
$ws.Range("Q10").Value = 3.1
$ws.Range("R10").Value = 1.36
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 8
$ws.Range("Z10").Value = 19
$ws.Range("AA10").Value = 23
$ws.Range("AC10").Value = 5
$ws.Range("AH10").Value = 8
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 17
$ws.Range("AL10").Value = 41
$ws.Range("AM10").Value = 51
$ws.Range("AN10").Value = 3.75
$ws.Range("AO10").Value = 13
$ws.Range("AP10").Value = 34
$ws.Range("AY10").Value = 26
$ws.Range("BA10").Value = 101
$ws.Range("BB10").Value = 151
